$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the hyperlinks that belong to the project rows (14-16) which are
#    about to be cleared out. Re-scan the live collection for each address
#    and delete one at a time (deleting while enumerating a stale snapshot
#    skips/misfires on this host, so we look the hyperlink up fresh each
#    time).
# ---------------------------------------------------------------------------
$staleLinks = @('$E$14', '$E$15', '$E$16', '$I$16')
foreach ($addr in $staleLinks) {
    $found = $null
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $found = $h
        }
    }
    if ($found -ne $null) {
        $found.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2. Remove stale "in progress" project rows (Dev Tool CLI / Mapping
#    Earthquakes / Biodiversity Dashboard) from rows 14-16, leaving behind
#    only the blank placeholder cells (with their original formatting) that
#    are still present in the template, exactly like the existing blank
#    row 10.
# ---------------------------------------------------------------------------
$ws.Range("A14:J16").ClearContents()

$ws.Range("A14").Clear()
$ws.Range("C14:D14").Clear()
$ws.Range("F14:H14").Clear()

$ws.Range("A15").Clear()
$ws.Range("C15:D15").Clear()
$ws.Range("F15:H15").Clear()

$ws.Range("A16").Clear()
$ws.Range("C16:D16").Clear()
$ws.Range("F16:H16").Clear()
$ws.Range("J16").Clear()

# ---------------------------------------------------------------------------
# 3. Update the concepts list for the "Setup/Workflow Guides" module (row 18)
#    now that the Dev Tool CLI project (and its docs) have been folded in.
# ---------------------------------------------------------------------------
$ws.Range("H18").Value = "documentation, project setup, workflows"

# ---------------------------------------------------------------------------
# 4. Restore the active selection to the last data row, matching where the
#    author left the cursor after editing.
# ---------------------------------------------------------------------------
$null = $ws.Range("A18:H18").Select()
